$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 19:33"

# --- Swap country labels (reorder within "País" list) ---
# Rows 8/9 : Colombia <-> España
$ws.Range("A8").Value = "España"
$ws.Range("A9").Value = "Colombia"

# Rows 73/74 : Kenia <-> Irlanda
$ws.Range("A73").Value = "Irlanda"
$ws.Range("A74").Value = "Kenia"

function Set-RowStats($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Row 4 : Estados Unidos
Set-RowStats 4 7858904 25141 5039542 2601252 0 372 218110

# Row 5 : India
Set-RowStats 5 6957511 53699 5955462 895010 0 518 107039

# Row 8 : now España
Set-RowStats 8 890367 5986 0 0 0 241 32929

# Row 9 : now Colombia
Set-RowStats 9 886179 0 777658 81190 0 0 27331

# Row 23
Set-RowStats 23 332382 1629 291754 31906 0 55 8722

# Row 65
Set-RowStats 65 52804 146 37067 13948 0 6 1789

# Row 73 : now Irlanda
Set-RowStats 73 40703 617 23364 15518 0 4 1821

# Row 74 : now Kenia
Set-RowStats 74 40178 0 31710 7717 0 0 751

# Row 109
Set-RowStats 109 9742 103 7171 2502 0 1 69

# Row 180 (only D and E changed)
$ws.Cells.Item(180, 4).Value = 459
$ws.Cells.Item(180, 5).Value = 18
